$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2053.6206
$ws.Range("I15").Value = 2053.6206
$ws.Range("K15").Value = 6160.861800000001
$ws.Range("M15").Value = -5991.861800000001
$ws.Range("H80").Value = 522.8125
$ws.Range("I80").Value = 376.1875
$ws.Range("J80").Value = 669.4375
$ws.Range("K80").Value = 1128.5625
$ws.Range("L80").Value = 2008.3125
$ws.Range("M80").Value = -130.5625
$ws.Range("N80").Value = -4004.3125
$ws.Range("H83").Value = 522.8125
$ws.Range("I83").Value = 376.1875
$ws.Range("J83").Value = 669.4375
$ws.Range("K83").Value = 3385.6875
$ws.Range("L83").Value = 6024.9375
$ws.Range("M83").Value = 1606.3125
$ws.Range("N83").Value = -16008.9375
$ws.Range("H86").Value = 3686.25
$ws.Range("I86").Value = 2491.1667
$ws.Range("K86").Value = 2491.1667
$ws.Range("M86").Value = -1368.1667
$ws.Range("H88").Value = 3383.2354
$ws.Range("J88").Value = 3727.4285
$ws.Range("L88").Value = 3727.4285
$ws.Range("N88").Value = -4539.4285
$ws.Range("H89").Value = 3686.25
$ws.Range("I89").Value = 2491.1667
$ws.Range("K89").Value = 12455.8335
$ws.Range("M89").Value = -6839.833500000001
$ws.Range("H91").Value = 3383.2354
$ws.Range("J91").Value = 3727.4285
$ws.Range("L91").Value = 3727.4285
$ws.Range("N91").Value = -6535.4285
$ws.Range("H112").Value = 3208.3572
$ws.Range("J112").Value = 3230.3455
$ws.Range("L112").Value = 9691.0365
$ws.Range("N112").Value = -11907.0365
$ws.Range("H113").Value = 3029.3333
$ws.Range("I113").Value = 2093.6
$ws.Range("J113").Value = 4199
$ws.Range("K113").Value = 2093.6
$ws.Range("L113").Value = 4199
$ws.Range("M113").Value = 1160.4
$ws.Range("N113").Value = -10707
$ws.Range("H138").Value = 2591.039
$ws.Range("I138").Value = 1472.7727
$ws.Range("J138").Value = 3038.3455
$ws.Range("K138").Value = 4418.3181
$ws.Range("L138").Value = 9115.0365
$ws.Range("M138").Value = 721.6818999999996
$ws.Range("N138").Value = -19395.0365

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1774.0555
$ws.Range("I2").Value = 1506.0714
$ws.Range("J2").Value = 2712
$ws.Range("K2").Value = 1506.0714
$ws.Range("L2").Value = 2712
$ws.Range("M2").Value = -1393.0714
$ws.Range("N2").Value = -2938
$ws.Range("H32").Value = 8544.718000000001
$ws.Range("I32").Value = 4077.3794
$ws.Range("J32").Value = 21500
$ws.Range("K32").Value = 4077.3794
$ws.Range("L32").Value = 21500
$ws.Range("M32").Value = -3790.3794
$ws.Range("N32").Value = -22074
$ws.Range("H45").Value = 3560.6965
$ws.Range("I45").Value = 2117.9167
$ws.Range("J45").Value = 3954.182
$ws.Range("K45").Value = 2117.9167
$ws.Range("L45").Value = 3954.182
$ws.Range("M45").Value = -1740.9167
$ws.Range("N45").Value = -4708.182
$ws.Range("H61").Value = 5577.0386
$ws.Range("I61").Value = 3727
$ws.Range("K61").Value = 3727
$ws.Range("M61").Value = -3515
$ws.Range("H88").Value = 2669.6365
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 3108.25
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 3108.25
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -3920.25
$ws.Range("H91").Value = 2669.6365
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 3108.25
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 3108.25
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -5916.25
$ws.Range("H116").Value = 1774.0555
$ws.Range("I116").Value = 1506.0714
$ws.Range("J116").Value = 2712
$ws.Range("K116").Value = 1506.0714
$ws.Range("L116").Value = 2712
$ws.Range("M116").Value = 787.9286
$ws.Range("N116").Value = -7300
$ws.Range("H136").Value = 5577.0386
$ws.Range("I136").Value = 3727
$ws.Range("K136").Value = 11181
$ws.Range("M136").Value = -8631

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1774.0555
$ws.Range("I3").Value = 1506.0714
$ws.Range("J3").Value = 2712
$ws.Range("K3").Value = 1506.0714
$ws.Range("L3").Value = 2712
$ws.Range("M3").Value = -1392.0714
$ws.Range("N3").Value = -2940
$ws.Range("H76").Value = 33093.145
$ws.Range("J76").Value = 33715.77
$ws.Range("L76").Value = 33715.77
$ws.Range("N76").Value = -34345.77
$ws.Range("H79").Value = 33093.145
$ws.Range("J79").Value = 33715.77
$ws.Range("L79").Value = 33715.77
$ws.Range("N79").Value = -35899.77
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H134").Value = 3748.4783
$ws.Range("I134").Value = 3629.2856
$ws.Range("K134").Value = 10887.8568
$ws.Range("M134").Value = -8352.856800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6094.222
$ws.Range("I62").Value = 6100.5
$ws.Range("K62").Value = 6100.5
$ws.Range("M62").Value = -5476.5
$ws.Range("H65").Value = 6094.222
$ws.Range("I65").Value = 6100.5
$ws.Range("K65").Value = 30502.5
$ws.Range("M65").Value = -27382.5
$ws.Range("H107").Value = 23810818
$ws.Range("I107").Value = 33334704
$ws.Range("J107").Value = 1103.6666
$ws.Range("K107").Value = 33334704
$ws.Range("L107").Value = 1103.6666
$ws.Range("M107").Value = -33332784
$ws.Range("N107").Value = -4943.6666
$ws.Range("H134").Value = 1696
$ws.Range("I134").Value = 1696.5
$ws.Range("J134").Value = 1694
$ws.Range("K134").Value = 5089.5
$ws.Range("L134").Value = 5082
$ws.Range("M134").Value = -2554.5
$ws.Range("N134").Value = -10152

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2012.0358
$ws.Range("I131").Value = 1191.4546
$ws.Range("J131").Value = 2543
$ws.Range("K131").Value = 3574.3638
$ws.Range("L131").Value = 7629
$ws.Range("M131").Value = 1465.6362
$ws.Range("N131").Value = -17709

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 31999
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H102").Value = 2036.8334
$ws.Range("I102").Value = 1791.5
$ws.Range("K102").Value = 1791.5
$ws.Range("M102").Value = -169.5
$ws.Range("H126").Value = 7357.5713
$ws.Range("I126").Value = 7103.625
$ws.Range("J126").Value = 7696.1665
$ws.Range("K126").Value = 21310.875
$ws.Range("L126").Value = 23088.4995
$ws.Range("M126").Value = -18840.875
$ws.Range("N126").Value = -28028.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 10749.75
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H14").Value = 793749.75
$ws.Range("J14").Value = 999999
$ws.Range("L14").Value = 999999
$ws.Range("N14").Value = -1000343
$ws.Range("H20").Value = 2793333.2
$ws.Range("H38").Value = 1000047
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 925
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 925
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 925
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1301
$ws.Range("H61").Value = 1043.5555
$ws.Range("I61").Value = 1074.375
$ws.Range("J61").Value = 797
$ws.Range("K61").Value = 1074.375
$ws.Range("L61").Value = 797
$ws.Range("M61").Value = -872.375
$ws.Range("N61").Value = -1201
$ws.Range("H100").Value = 8000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 8000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 8000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -9082
$ws.Range("H113").Value = 1043.5555
$ws.Range("I113").Value = 1074.375
$ws.Range("J113").Value = 797
$ws.Range("K113").Value = 1074.375
$ws.Range("L113").Value = 797
$ws.Range("M113").Value = 1095.625
$ws.Range("N113").Value = -5137
$ws.Range("H133").Value = 89994.60000000001
$ws.Range("J133").Value = 89994.60000000001
$ws.Range("L133").Value = 89994.60000000001
$ws.Range("N133").Value = -95054.60000000001
$ws.Range("H136").Value = 3885.8696
$ws.Range("I136").Value = 4159.222
$ws.Range("K136").Value = 12477.666
$ws.Range("M136").Value = -9927.665999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 1454.5
$ws.Range("I10").Value = 1454.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1454.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -1285.5
$ws.Range("N10").ClearContents()
$ws.Range("H122").Value = 4548.231
$ws.Range("I122").Value = 3248.7778
$ws.Range("K122").Value = 9746.3334
$ws.Range("M122").Value = -7296.3334
$ws.Range("H136").Value = 3293.842
$ws.Range("I136").Value = 3147.8064
$ws.Range("K136").Value = 9443.4192
$ws.Range("M136").Value = -6893.4192
